$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.345.74"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.22%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'1.865.44"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.31%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("E4").Value = "'  +0.22%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'330.45"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -2.10%  "
$ws.Range("E5").Style = "Normal"

$ws.Range("E6").Value = "'  +0.12%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("E7").Value = "'  -2.64%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = "'0.4005"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +1.82%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = "'47.67"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +0.91%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Value = "'0.07837"
$ws.Range("D10").Style = "Normal"

$ws.Range("D11").Value = "'0.9822"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -2.22%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = "'21.18"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -2.69%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = "'1.873.86"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +0.02%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'5.824"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -2.81%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'6.988"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -4.15%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("E16").Value = "'  +0.16%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Value = "'88.12"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -3.45%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = "'0.06544"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -0.58%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("E19").Value = "'  -2.49%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").Value = "'17.14"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -3.04%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = "'1.003"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +0.19%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Value = "'28.332.33"
$ws.Range("D22").Style = "Normal"

$ws.Range("D23").Value = "'5.326"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -2.28%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = "'10.82"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -2.00%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = "'2.251"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -1.97%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").Value = "'2.099.15"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +0.07%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("D27").Value = "'157.24"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -1.31%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("E28").Value = "'  -2.82%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("D29").Value = "'2.056"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -4.73%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("D30").Value = "'5.283"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -3.99%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("D31").Value = "'117.04"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -2.44%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("D32").Value = "'0.9518"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -2.90%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("D33").Value = "'0.09310"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -1.97%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("D34").Value = "'3.596"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +0.36%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("D35").Value = "'1.380"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -0.11%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("D36").Value = "'5.215"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -2.74%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("D37").Value = "'0.06015"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -1.23%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("D38").Value = "'0.02199"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -3.35%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("D39").Value = "'8.271"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -2.15%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").Value = "'1.162"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -1.28%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").Value = "'1.002"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +0.15%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("D42").Value = "'0.5734"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -3.93%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").Value = "'0.1803"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -3.96%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = "'9.996"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -3.90%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = "'1.259"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -3.42%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").Value = "'2.276"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +12.70%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("D47").Value = "'0.5409"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -3.75%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").Value = "'11.85"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -2.49%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("D49").Value = "'0.07171"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +3.96%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("D50").Value = "'1.877"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -4.70%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("D51").Value = "'109.92"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -0.73%  "
$ws.Range("E51").Style = "Normal"
